$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36 - this shifts the existing rows 36-73 down to 37-74,
# preserving all of their data/formatting (including the date-formatted column D style).
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Range("A36").Value = 6
$ws.Range("B36").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 44778
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = "Tropicales y subtropicales"
$ws.Range("I36").Value = 100108007
$ws.Range("J36").Value = "Coco"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 150
$ws.Range("N36").Value = 24000
$ws.Range("O36").Value = 24000
$ws.Range("P36").Value = 24000
$ws.Range("Q36").Value = "$/malla 20 unidades"
$ws.Range("R36").Value = "Perú"
$ws.Range("S36").Value = 1200
$ws.Range("T36").Value = 20
